$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Tgfb1"
$row2[0,2] = "Itgb6"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 51.59157666666666
$row2[0,7] = 154.77473
$row2[0,8] = 0.2641250550177587
$row2[0,9] = 0.2641250550177588
$row2[0,10] = 2
$row2[0,11] = 0.6666666666666666
$row2[0,12] = 0.3530253333333334
$row2[0,13] = 1.059076
$row2[0,14] = 0.05211323950342944
$row2[0,15] = 0.05211323950342944
$row2[0,16] = 18.21313354994222
$row2[0,17] = 163.91820194948
$row2[0,18] = 0.01376441225099694
$row2[0,19] = 0.01376441225099694
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Tgfb1"
$row3[0,2] = "Itgb6"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 51.59157666666666
$row3[0,7] = 154.77473
$row3[0,8] = 0.2641250550177587
$row3[0,9] = 0.2641250550177588
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 2.310677666666666
$row3[0,13] = 6.932033
$row3[0,14] = 0.3410998794937062
$row3[0,15] = 0.3410998794937063
$row3[0,16] = 119.2115039917878
$row3[0,17] = 1072.90353592609
$row3[0,18] = 0.09009302443782603
$row3[0,19] = 0.09009302443782607
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Tgfb1"
$row4[0,2] = "Itgb6"
$row4[0,3] = "M2"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 51.59157666666666
$row4[0,7] = 154.77473
$row4[0,8] = 0.2641250550177587
$row4[0,9] = 0.2641250550177588
$row4[0,10] = 1
$row4[0,11] = 0.3333333333333333
$row4[0,12] = 0.07663533333333333
$row4[0,13] = 0.229906
$row4[0,14] = 0.01131282971314188
$row4[0,15] = 0.01131282971314188
$row4[0,16] = 3.953737675042222
$row4[0,17] = 35.58363907537999
$row4[0,18] = 0.002988001770390135
$row4[0,19] = 0.002988001770390135
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Tgfb1"
$row5[0,2] = "Itgb6"
$row5[0,3] = "sCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 51.59157666666666
$row5[0,7] = 154.77473
$row5[0,8] = 0.2641250550177587
$row5[0,9] = 0.2641250550177588
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 4.033858333333334
$row5[0,13] = 12.101575
$row5[0,14] = 0.5954740512897225
$row5[0,15] = 0.5954740512897225
$row5[0,16] = 208.1131114666389
$row5[0,17] = 1873.01800319975
$row5[0,18] = 0.1572796165585456
$row5[0,19] = 0.1572796165585457
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Tgfb1"
$row6[0,2] = "Itgb6"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 19.32115333333334
$row6[0,7] = 57.96346000000001
$row6[0,8] = 0.09891538535728452
$row6[0,9] = 0.09891538535728453
$row6[0,10] = 2
$row6[0,11] = 0.6666666666666666
$row6[0,12] = 0.3530253333333334
$row6[0,13] = 1.059076
$row6[0,14] = 0.05211323950342944
$row6[0,15] = 0.05211323950342944
$row6[0,16] = 6.820856595884447
$row6[0,17] = 61.38770936296002
$row6[0,18] = 0.005154801167698185
$row6[0,19] = 0.005154801167698187
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Tgfb1"
$row7[0,2] = "Itgb6"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 19.32115333333334
$row7[0,7] = 57.96346000000001
$row7[0,8] = 0.09891538535728452
$row7[0,9] = 0.09891538535728453
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 2.310677666666666
$row7[0,13] = 6.932033
$row7[0,14] = 0.3410998794937062
$row7[0,15] = 0.3410998794937063
$row7[0,16] = 44.64495750157556
$row7[0,17] = 401.8046175141801
$row7[0,18] = 0.03374002602544326
$row7[0,19] = 0.03374002602544327
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Tgfb1"
$row8[0,2] = "Itgb6"
$row8[0,3] = "M2"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 19.32115333333334
$row8[0,7] = 57.96346000000001
$row8[0,8] = 0.09891538535728452
$row8[0,9] = 0.09891538535728453
$row8[0,10] = 1
$row8[0,11] = 0.3333333333333333
$row8[0,12] = 0.07663533333333333
$row8[0,13] = 0.229906
$row8[0,14] = 0.01131282971314188
$row8[0,15] = 0.01131282971314188
$row8[0,16] = 1.480683026084445
$row8[0,17] = 13.32614723476
$row8[0,18] = 0.001119012910556767
$row8[0,19] = 0.001119012910556768
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Tgfb1"
$row9[0,2] = "Itgb6"
$row9[0,3] = "sCs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 19.32115333333334
$row9[0,7] = 57.96346000000001
$row9[0,8] = 0.09891538535728452
$row9[0,9] = 0.09891538535728453
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 4.033858333333334
$row9[0,13] = 12.101575
$row9[0,14] = 0.5954740512897225
$row9[0,15] = 0.5954740512897225
$row9[0,16] = 77.9387953832778
$row9[0,17] = 701.4491584495001
$row9[0,18] = 0.05890154525358631
$row9[0,19] = 0.05890154525358632
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "M2"
$row10[0,1] = "Tgfb1"
$row10[0,2] = "Itgb6"
$row10[0,3] = "ECs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 112.3724673333333
$row10[0,7] = 337.117402
$row10[0,8] = 0.5752951554216499
$row10[0,9] = 0.57529515542165
$row10[0,10] = 2
$row10[0,11] = 0.6666666666666666
$row10[0,12] = 0.3530253333333334
$row10[0,13] = 1.059076
$row10[0,14] = 0.05211323950342944
$row10[0,15] = 0.05211323950342944
$row10[0,16] = 39.67032773783911
$row10[0,17] = 357.032949640552
$row10[0,18] = 0.0299804942196511
$row10[0,19] = 0.02998049421965111
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "M2"
$row11[0,1] = "Tgfb1"
$row11[0,2] = "Itgb6"
$row11[0,3] = "FAPs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 112.3724673333333
$row11[0,7] = 337.117402
$row11[0,8] = 0.5752951554216499
$row11[0,9] = 0.57529515542165
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 2.310677666666666
$row11[0,13] = 6.932033
$row11[0,14] = 0.3410998794937062
$row11[0,15] = 0.3410998794937063
$row11[0,16] = 259.6565506153628
$row11[0,17] = 2336.908955538266
$row11[0,18] = 0.1962331081876378
$row11[0,19] = 0.1962331081876378
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "M2"
$row12[0,1] = "Tgfb1"
$row12[0,2] = "Itgb6"
$row12[0,3] = "M2"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 112.3724673333333
$row12[0,7] = 337.117402
$row12[0,8] = 0.5752951554216499
$row12[0,9] = 0.57529515542165
$row12[0,10] = 1
$row12[0,11] = 0.3333333333333333
$row12[0,12] = 0.07663533333333333
$row12[0,13] = 0.229906
$row12[0,14] = 0.01131282971314188
$row12[0,15] = 0.01131282971314188
$row12[0,16] = 8.61170149157911
$row12[0,17] = 77.50531342421199
$row12[0,18] = 0.006508216128080616
$row12[0,19] = 0.006508216128080617
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "M2"
$row13[0,1] = "Tgfb1"
$row13[0,2] = "Itgb6"
$row13[0,3] = "sCs"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 112.3724673333333
$row13[0,7] = 337.117402
$row13[0,8] = 0.5752951554216499
$row13[0,9] = 0.57529515542165
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 4.033858333333334
$row13[0,13] = 12.101575
$row13[0,14] = 0.5954740512897225
$row13[0,15] = 0.5954740512897225
$row13[0,16] = 453.2946137897944
$row13[0,17] = 4079.65152410815
$row13[0,18] = 0.3425733368862804
$row13[0,19] = 0.3425733368862805
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,20
$row14[0,0] = "sCs"
$row14[0,1] = "Tgfb1"
$row14[0,2] = "Itgb6"
$row14[0,3] = "ECs"
$row14[0,4] = 3
$row14[0,5] = 1
$row14[0,6] = 12.044915
$row14[0,7] = 36.134745
$row14[0,8] = 0.06166440420330686
$row14[0,9] = 0.06166440420330688
$row14[0,10] = 2
$row14[0,11] = 0.6666666666666666
$row14[0,12] = 0.3530253333333334
$row14[0,13] = 1.059076
$row14[0,14] = 0.05211323950342944
$row14[0,15] = 0.05211323950342944
$row14[0,16] = 4.252160132846668
$row14[0,17] = 38.26944119562
$row14[0,18] = 0.003213531865083211
$row14[0,19] = 0.003213531865083212
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,20
$row15[0,0] = "sCs"
$row15[0,1] = "Tgfb1"
$row15[0,2] = "Itgb6"
$row15[0,3] = "FAPs"
$row15[0,4] = 3
$row15[0,5] = 1
$row15[0,6] = 12.044915
$row15[0,7] = 36.134745
$row15[0,8] = 0.06166440420330686
$row15[0,9] = 0.06166440420330688
$row15[0,10] = 3
$row15[0,11] = 1
$row15[0,12] = 2.310677666666666
$row15[0,13] = 6.932033
$row15[0,14] = 0.3410998794937062
$row15[0,15] = 0.3410998794937063
$row15[0,16] = 27.83191608739833
$row15[0,17] = 250.487244786585
$row15[0,18] = 0.02103372084279916
$row15[0,19] = 0.02103372084279917
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,20
$row16[0,0] = "sCs"
$row16[0,1] = "Tgfb1"
$row16[0,2] = "Itgb6"
$row16[0,3] = "M2"
$row16[0,4] = 3
$row16[0,5] = 1
$row16[0,6] = 12.044915
$row16[0,7] = 36.134745
$row16[0,8] = 0.06166440420330686
$row16[0,9] = 0.06166440420330688
$row16[0,10] = 1
$row16[0,11] = 0.3333333333333333
$row16[0,12] = 0.07663533333333333
$row16[0,13] = 0.229906
$row16[0,14] = 0.01131282971314188
$row16[0,15] = 0.01131282971314188
$row16[0,16] = 0.9230660759966668
$row16[0,17] = 8.30759468397
$row16[0,18] = 0.0006975989041143609
$row16[0,19] = 0.000697598904114361
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,20
$row17[0,0] = "sCs"
$row17[0,1] = "Tgfb1"
$row17[0,2] = "Itgb6"
$row17[0,3] = "sCs"
$row17[0,4] = 3
$row17[0,5] = 1
$row17[0,6] = 12.044915
$row17[0,7] = 36.134745
$row17[0,8] = 0.06166440420330686
$row17[0,9] = 0.06166440420330688
$row17[0,10] = 3
$row17[0,11] = 1
$row17[0,12] = 4.033858333333334
$row17[0,13] = 12.101575
$row17[0,14] = 0.5954740512897225
$row17[0,15] = 0.5954740512897225
$row17[0,16] = 48.58748074704167
$row17[0,17] = 437.287326723375
$row17[0,18] = 0.03671955259131014
$row17[0,19] = 0.03671955259131014
$ws.Range("A17:T17").Value = $row17
